$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 580; this shifts the old rows 580-603 down to 582-605
$ws.Rows.Item(580).Resize(2).Insert()

# Populate the two new rows (580, 581) with the new weekly price entries
# Row 580
$ws.Cells.Item(580, 1).Value2 = 11
$ws.Cells.Item(580, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(580, 3).Value2 = 'Bíobío'
$ws.Cells.Item(580, 4).Value2 = 44939
$ws.Cells.Item(580, 5).Value2 = 8
$ws.Cells.Item(580, 6).Value2 = 100112020
$ws.Cells.Item(580, 7).Value2 = 'Tomate'
$ws.Cells.Item(580, 8).Value2 = 'Semiduro'
$ws.Cells.Item(580, 9).Value2 = 'Primera'
$ws.Cells.Item(580, 10).Value2 = 600
$ws.Cells.Item(580, 11).Value2 = 9000
$ws.Cells.Item(580, 12).Value2 = 10000
$ws.Cells.Item(580, 13).Value2 = 9500
$ws.Cells.Item(580, 14).Value2 = '$/bandeja 18 kilos'
$ws.Cells.Item(580, 15).Value2 = 'Región del Maule'
$ws.Cells.Item(580, 16).Value2 = 528
$ws.Cells.Item(580, 17).Value2 = 18
$ws.Cells.Item(580, 18).Value2 = 'Hortaliza'

# Row 581
$ws.Cells.Item(581, 1).Value2 = 11
$ws.Cells.Item(581, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(581, 3).Value2 = 'Bíobío'
$ws.Cells.Item(581, 4).Value2 = 44939
$ws.Cells.Item(581, 5).Value2 = 8
$ws.Cells.Item(581, 6).Value2 = 100112020
$ws.Cells.Item(581, 7).Value2 = 'Tomate'
$ws.Cells.Item(581, 8).Value2 = 'Semiduro'
$ws.Cells.Item(581, 9).Value2 = 'Segunda'
$ws.Cells.Item(581, 10).Value2 = 300
$ws.Cells.Item(581, 11).Value2 = 8000
$ws.Cells.Item(581, 12).Value2 = 8000
$ws.Cells.Item(581, 13).Value2 = 8000
$ws.Cells.Item(581, 14).Value2 = '$/bandeja 18 kilos'
$ws.Cells.Item(581, 15).Value2 = 'Región del Maule'
$ws.Cells.Item(581, 16).Value2 = 444
$ws.Cells.Item(581, 17).Value2 = 18
$ws.Cells.Item(581, 18).Value2 = 'Hortaliza'

